# 2020 annual report update
# Update the 2019 "CAP Catch/hr*" values for Center, North Twin, South Twin
# and Five Island lakes (column C, rows 14-17), then move the active
# selection to reflect where the editor was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C14").Value = 2.86
$ws.Range("C15").Value = 8.17
$ws.Range("C16").Value = 5.33
$ws.Range("C17").Value = 32.5

[void]$ws.Range("C6").Select()
